# ValueSet-med-defibrotide-vs.xlsx update: new logo and colors
# (commit message is generic; the actual edit refreshes the IG-publisher
#  generated "Metadata" sheet: version, status, date, and contact info,
#  including a new Jurisdiction row - see the authoritative xml diff.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update existing metadata values --------------------------------------
$ws.Range("B3").Value  = "0.1.7"                                    # Version
$ws.Range("B6").Value  = "draft"                                    # Status
$ws.Range("B8").Value  = "2024-08-23T10:17:11-05:00"                # Date
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"  # Contact

# Row 11 used to be a duplicate "Contact" row with the same value as row 10;
# turn it into the new second Contact entry (Bob Milius).
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row (row 12), empty value -----------------
$ws.Rows.Item(12).Insert()
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$excel.CutCopyMode = $false
